# Updates the cryptos price/volume table to the latest scraped values.
# Leading "'" forces Excel to store these as text (matching the original
# inlineStr/text cells) instead of re-interpreting numeric-looking
# strings (e.g. "1.010", "0.4700") as numbers and losing trailing zeros
# or the multi-dot "thousands" formatted prices (e.g. "27.646.89").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.646.89"
$ws.Range("E2").Value = "'  -0.27%  "
$ws.Range("D3").Value = "'1.864.21"
$ws.Range("E3").Value = "'  -0.84%  "
$ws.Range("D4").Value = "'1.010"
$ws.Range("E4").Value = "'  +0.24%  "
$ws.Range("D5").Value = "'333.39"
$ws.Range("E5").Value = "'  -0.01%  "
$ws.Range("D6").Value = "'1.009"
$ws.Range("E6").Value = "'  +0.26%  "
$ws.Range("D7").Value = "'0.4700"
$ws.Range("E7").Value = "'  -0.13%  "
$ws.Range("D8").Value = "'0.3921"
$ws.Range("E8").Value = "'  -0.41%  "
$ws.Range("B9").Value = "'Dogecoin"
$ws.Range("C9").Value = "'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D9").Value = "'0.07990"
$ws.Range("E9").Value = "'  -0.87%  "
$ws.Range("B10").Value = "'OKB"
$ws.Range("C10").Value = "'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").Value = "'45.27"
$ws.Range("E10").Value = "'  -4.74%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "'  -2.44%  "
$ws.Range("D12").Value = "'21.79"
$ws.Range("E12").Value = "'  -1.75%  "
$ws.Range("D13").Value = "'5.995"
$ws.Range("E13").Value = "'  +0.32%  "
$ws.Range("D14").Value = "'1.857.29"
$ws.Range("E14").Value = "'  -1.59%  "
$ws.Range("E15").Value = "'  +1.49%  "
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("D17").Value = "'88.51"
$ws.Range("E17").Value = "'  +1.63%  "
$ws.Range("D18").Value = "'0.06715"
$ws.Range("E18").Value = "'  +0.14%  "
$ws.Range("E19").Value = "'  -0.73%  "
$ws.Range("D20").Value = "'17.12"
$ws.Range("E20").Value = "'  -0.73%  "
$ws.Range("D21").Value = "'1.010"
$ws.Range("E21").Value = "'  +0.27%  "
$ws.Range("D22").Value = "'27.630.33"
$ws.Range("E22").Value = "'  -0.38%  "
$ws.Range("E23").Value = "'  -1.55%  "
$ws.Range("E24").Value = "'  -1.19%  "
$ws.Range("D25").Value = "'2.314"
$ws.Range("E25").Value = "'  -0.10%  "
$ws.Range("D26").Value = "'2.082.28"
$ws.Range("E26").Value = "'  -1.37%  "
$ws.Range("E27").Value = "'  -0.66%  "
$ws.Range("E28").Value = "'  -2.05%  "
$ws.Range("D29").Value = "'2.148"
$ws.Range("E29").Value = "'  +1.95%  "
$ws.Range("D30").Value = "'5.434"
$ws.Range("E30").Value = "'  -2.87%  "
$ws.Range("E31").Value = "'  -0.14%  "
$ws.Range("D32").Value = "'0.9818"
$ws.Range("E32").Value = "'  -0.33%  "
$ws.Range("D33").Value = "'0.09488"
$ws.Range("E33").Value = "'  +0.09%  "
$ws.Range("D34").Value = "'3.615"
$ws.Range("E34").Value = "'  +0.00%  "
$ws.Range("D35").Value = "'5.305"
$ws.Range("E35").Value = "'  -1.02%  "
$ws.Range("D36").Value = "'1.336"
$ws.Range("E36").Value = "'  -8.16%  "
$ws.Range("D37").Value = "'0.06052"
$ws.Range("E37").Value = "'  -1.49%  "
$ws.Range("D38").Value = "'0.02229"
$ws.Range("E38").Value = "'  -1.62%  "
$ws.Range("D39").Value = "'8.319"
$ws.Range("E39").Value = "'  +2.34%  "
$ws.Range("D40").Value = "'1.195"
$ws.Range("E40").Value = "'  -2.92%  "
$ws.Range("D41").Value = "'1.010"
$ws.Range("E41").Value = "'  +0.30%  "
$ws.Range("D42").Value = "'0.5963"
$ws.Range("E42").Value = "'  -0.68%  "
$ws.Range("D43").Value = "'0.1885"
$ws.Range("E43").Value = "'  -0.82%  "
$ws.Range("E44").Value = "'  +0.05%  "
$ws.Range("D45").Value = "'1.249"
$ws.Range("E45").Value = "'  -0.64%  "
$ws.Range("D46").Value = "'0.5645"
$ws.Range("E46").Value = "'  -1.04%  "
$ws.Range("D47").Value = "'12.19"
$ws.Range("E47").Value = "'  +0.12%  "
$ws.Range("D48").Value = "'1.921"
$ws.Range("E48").Value = "'  -1.41%  "
$ws.Range("D49").Value = "'0.06755"
$ws.Range("E49").Value = "'  -2.14%  "
$ws.Range("D50").Value = "'111.83"
$ws.Range("E50").Value = "'  -2.37%  "
$ws.Range("D51").Value = "'3.082"
$ws.Range("E51").Value = "'  -9.44%  "
